$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 59: "Outwitting the Devil" by Napolean Hill (Self-Development) ---
$ws.Range("E58:H58").Copy($ws.Range("E59:H59"))
$ws.Cells.Item(59, 5).Value = 55
$ws.Cells.Item(59, 6).Value = "Outwitting the Devil"
$ws.Cells.Item(59, 7).Value = "Napolean Hill"
$ws.Cells.Item(59, 8).Value = "Self-Development"
$ws.Rows("59").RowHeight = 21

# --- Row 60: "The Oaths of Vayuputras - Shiva trilogy book-3" by Amish (Mythology) ---
$ws.Range("E58:H58").Copy($ws.Range("E60:H60"))
$ws.Cells.Item(60, 5).Value = 56
$ws.Cells.Item(60, 6).Value = "The Oaths of Vayuputras - Shiva trilogy book-3"
$ws.Cells.Item(60, 7).Value = "Amish "
$ws.Cells.Item(60, 8).Value = "Mythology"
$ws.Rows("60").RowHeight = 21

# --- Row 61: trailing empty templated row (same fills, borders stripped) ---
$ws.Range("E58:H58").Copy($ws.Range("E61:H61"))
$ws.Range("E61:H61").ClearContents()
$ws.Range("E61:H61").Borders.LineStyle = 0
$ws.Rows("61").RowHeight = 21

# --- Update the selection to match where the user ended up after typing ---
$ws.Range("G66").Select() | Out-Null
